# Apply the cryptos-list refresh described by the commit diff.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). All cell contents are
# plain text in the source workbook (prices/links stored as strings, not
# numbers), so every write below goes through Set-TextCell, which forces
# Excel to keep the value as text (a leading apostrophe marks it a text
# entry) and then resets the cell style back to Normal so no stray
# "quote prefix" / number-format styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $cell = $ws.Range($addr)
    if ($text -match '^-?\d+(\.\d+)?$') {
        # Looks like a plain number (e.g. "6.60", "0.120") -- without help
        # Excel would parse this as a numeric value and silently drop
        # formatting such as trailing zeros. Force text entry instead.
        $cell.Value = "'" + $text
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $text
    }
}

Set-TextCell 'D2' '62.821.86'
Set-TextCell 'E2' '  -1.07%  '
Set-TextCell 'D3' '3.084.37'
Set-TextCell 'E3' '  +1.01%  '
Set-TextCell 'E4' '  -0.14%  '
Set-TextCell 'D5' '551.81'
Set-TextCell 'E5' '  +0.60%  '
Set-TextCell 'D6' '136.88'
Set-TextCell 'E6' '  -2.47%  '
Set-TextCell 'D7' '0.999'
Set-TextCell 'E7' '  -0.14%  '
Set-TextCell 'D8' '3.074.40'
Set-TextCell 'E8' '  +0.89%  '
Set-TextCell 'D9' '0.492'
Set-TextCell 'E9' '  +1.17%  '
Set-TextCell 'D10' '6.60'
Set-TextCell 'E10' '  +2.29%  '
Set-TextCell 'E11' '  +5.69%  '
Set-TextCell 'D12' '0.452'
Set-TextCell 'E12' '  +1.92%  '
Set-TextCell 'D13' '34.85'
Set-TextCell 'E13' '  -1.25%  '
Set-TextCell 'D14' '0.0000216'
Set-TextCell 'E14' '  +1.67%  '
Set-TextCell 'D15' '3.576.39'
Set-TextCell 'D16' '62.860.94'
Set-TextCell 'E16' '  -1.15%  '
Set-TextCell 'E17' '  +0.25%  '
Set-TextCell 'D18' '3.080.00'
Set-TextCell 'E18' '  +0.70%  '
Set-TextCell 'D19' '500.66'
Set-TextCell 'E19' '  +3.22%  '
Set-TextCell 'D20' '6.63'
Set-TextCell 'E20' '  +1.90%  '
Set-TextCell 'D21' '13.45'
Set-TextCell 'E21' '  +0.11%  '
Set-TextCell 'D22' '0.703'
Set-TextCell 'E22' '  +4.46%  '
Set-TextCell 'D23' '7.21'
Set-TextCell 'E23' '  +1.41%  '
Set-TextCell 'E24' '  +0.58%  '
Set-TextCell 'D25' '12.22'
Set-TextCell 'E25' '  +0.02%  '
Set-TextCell 'D26' '0.999'
Set-TextCell 'E26' '  -0.08%  '
Set-TextCell 'E27' '  +2.37%  '
Set-TextCell 'D28' '8.12'
Set-TextCell 'E28' '  +0.60%  '
Set-TextCell 'D29' '0.999'
Set-TextCell 'E29' '  -0.19%  '
Set-TextCell 'E30' '  -4.06%  '
Set-TextCell 'D31' '26.13'
Set-TextCell 'E31' '  +2.44%  '
Set-TextCell 'B32' 'Stacks'
Set-TextCell 'C32' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 'D32' '2.49'
Set-TextCell 'E32' '  -3.37%  '
Set-TextCell 'B33' 'Mantle'
Set-TextCell 'C33' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D33' '1.11'
Set-TextCell 'E33' '  -0.15%  '
Set-TextCell 'D34' '58.90'
Set-TextCell 'E34' '  +13.45%  '
Set-TextCell 'D35' '527.81'
Set-TextCell 'E35' '  -7.82%  '
Set-TextCell 'D36' '5.87'
Set-TextCell 'E36' '  +1.34%  '
Set-TextCell 'D37' '5.16'
Set-TextCell 'E37' '  -1.83%  '
Set-TextCell 'D38' '0.0410'
Set-TextCell 'E38' '  +3.47%  '
Set-TextCell 'D39' '3.036.64'
Set-TextCell 'E39' '  +1.67%  '
Set-TextCell 'B40' 'Hedera'
Set-TextCell 'C40' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D40' '0.0784'
Set-TextCell 'E40' '  +0.42%  '
Set-TextCell 'B41' 'Kaspa'
Set-TextCell 'C41' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 'D41' '0.120'
Set-TextCell 'E41' '  +3.94%  '
Set-TextCell 'D42' '8.04'
Set-TextCell 'E42' '  +0.10%  '
Set-TextCell 'E43' '  -5.49%  '
Set-TextCell 'D44' '0.252'
Set-TextCell 'E44' '  +5.25%  '
Set-TextCell 'D46' '2.05'
Set-TextCell 'E46' '  +0.83%  '
Set-TextCell 'D47' '121.63'
Set-TextCell 'E47' '  +4.31%  '
Set-TextCell 'E48' '  -0.28%  '
Set-TextCell 'D49' '23.65'
Set-TextCell 'E49' '  -3.07%  '
Set-TextCell 'D50' '0.0₃0502'
Set-TextCell 'E50' '  -1.80%  '
Set-TextCell 'D51' '2.37'
Set-TextCell 'E51' '  +70.94%  '
